# Applies the odds updates described in the commit diff.
# Workbook has a single sheet ("Sheet1"); update specific cells in
# rows 5, 6, 7, 8, 9, 15, 47, 61, 62, 63.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Range("Q5").Value = 1.07
$ws.Range("R5").Value = 7.5
$ws.Range("U5").Value = 1.93
$ws.Range("V5").Value = 1.83

# Row 6
$ws.Range("Q6").Value = 3.4
$ws.Range("R6").Value = 1.33
$ws.Range("S6").Value = 1.73
$ws.Range("T6").Value = 2
$ws.Range("AC6").Value = 5
$ws.Range("AH6").Value = 7.5
$ws.Range("AT6").Value = 2
$ws.Range("AY6").Value = 26

# Row 7
$ws.Range("V7").Value = 1.57
$ws.Range("AW7").Value = 151

# Row 8
$ws.Range("K8").Value = 1.92
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("V8").Value = 1.57

# Row 9
$ws.Range("K9").Value = 1.77
$ws.Range("M9").Value = 1.14
$ws.Range("N9").Value = 5.5
$ws.Range("V9").Value = 1.5

# Row 15
$ws.Range("G15").Value = 1.45

# Row 47
$ws.Range("G47").Value = 1.27
$ws.Range("I47").Value = 9.25
$ws.Range("J47").Value = 1.7
$ws.Range("K47").Value = 2.5
$ws.Range("L47").Value = 7.6
$ws.Range("Q47").Value = 1.55
$ws.Range("R47").Value = 2.15
$ws.Range("T47").Value = 3.32
$ws.Range("U47").Value = 1.98
$ws.Range("V47").Value = 1.65
$ws.Range("W47").Value = 7.4
$ws.Range("X47").Value = 6.2
$ws.Range("AA47").Value = 10.75
$ws.Range("AB47").Value = 29
$ws.Range("AC47").Value = 14
$ws.Range("AF47").Value = 110
$ws.Range("AG47").Value = 900
$ws.Range("AH47").Value = 25
$ws.Range("AI47").Value = 70
$ws.Range("AL47").Value = 120
$ws.Range("AP47").Value = 16
$ws.Range("AR47").Value = 40
$ws.Range("AT47").Value = 3.15
$ws.Range("AU47").Value = 9
$ws.Range("AX47").Value = 9.5
$ws.Range("AY47").Value = 55

# Row 61
$ws.Range("R61").Value = 1.57

# Row 62
$ws.Range("R62").Value = 1.53

# Row 63
$ws.Range("R63").Value = 1.7
